$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the marksheet's "Marking" row and "Total" row correct-mark counts,
# and the corresponding "Correct/Total" summary text.
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 130
$ws.Range("E12").Value = "130/140"
